{"js": "// Resume content refresh: replace the two \"Angular 4 / AngularJS\" mentions\n// in the SUMMARY section with \"Angular and React\" wording.\n\nconst body = context.document.body;\n\n// 1) \"Full stack developer ...\" bullet.\nconst search1 = body.search(\n  \"Full stack developer with over 10 years of Enterprise development experience with 5+ years experience working with Node and Angular up to Angular 4. \",\n  { matchCase: true, matchWholeWord: false }\n);\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\n    \"Full stack developer with over 10 years of Enterprise development experience with 5+ years experience working with Node, Angular and React. \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) \"Experienced with AngularJS 1.x up through Angular 4, api design ...\" bullet.\nconst search2 = body.search(\n  \"Experienced with AngularJS 1.x up through Angular 4, api design and build with Node and Express, Bootstrap up to Bootstrap 4, MongoDB and various build libraries (Webpack, Grunt, Gulp). I am experienced with AWS having architected, built and managed various projects.\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\n    \"Experienced with Angular and React, api design and build with Node and Express, Bootstrap up to Bootstrap 4, MongoDB and various build libraries (Webpack, Grunt, Gulp). I am experienced with AWS having architected, built and managed various projects.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Resume content refresh: replace the two \"Angular 4 / AngularJS\" mentions\n# in the SUMMARY section with \"Angular and React\" wording.\n\n$d = $word.ActiveDocument\n\n# 1) \"Full stack developer ...\" bullet.\n$oldText1 = \"Full stack developer with over 10 years of Enterprise development experience with 5+ years experience working with Node and Angular up to Angular 4. \"\n$newText1 = \"Full stack developer with over 10 years of Enterprise development experience with 5+ years experience working with Node, Angular and React. \"\n$range1 = $d.Content\n$range1.Find.Execute($oldText1, $false, $false, $false, $false, $false, $true, 1, $false, $newText1, 2) | Out-Null\n\n# 2) \"Experienced with AngularJS 1.x up through Angular 4, api design ...\" bullet.\n$oldText2 = \"Experienced with AngularJS 1.x up through Angular 4, api design and build with Node and Express, Bootstrap up to Bootstrap 4, MongoDB and various build libraries (Webpack, Grunt, Gulp). I am experienced with AWS having architected, built and managed various projects.\"\n$newText2 = \"Experienced with Angular and React, api design and build with Node and Express, Bootstrap up to Bootstrap 4, MongoDB and various build libraries (Webpack, Grunt, Gulp). I am experienced with AWS having architected, built and managed various projects.\"\n$range2 = $d.Content\n$range2.Find.Execute($oldText2, $false, $false, $false, $false, $false, $true, 1, $false, $newText2, 2) | Out-Null\n\nWrite-Output \"done\"\n"}
